$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged content, kept for clarity)
$ws.Range("A1").Value = "Responsibility"
$ws.Range("B1").Value = "Person"
$ws.Range("C1").Value = "Room Number"

# Updated table of responsibilities
$ws.Range("A2").Value = "Couch & Speaker Area"
$ws.Range("B2").Value = "Robin Epple"
$ws.Range("C2").Value = "U 13"

$ws.Range("A3").Value = "Dishes & Washing Utilities"
$ws.Range("B3").Value = "Silvy Kurzendorfer"
$ws.Range("C3").Value = "U 17"

$ws.Range("A4").Value = "Floor & Window Sill"
$ws.Range("B4").Value = "Robin Epple"
$ws.Range("C4").Value = "U 13"

$ws.Range("A5").Value = "Shelves"

$ws.Range("A6").Value = "Freezer"
$ws.Range("B6").Value = "Luke Caputo G."
$ws.Range("C6").Value = "U 09"

$ws.Range("A7").Value = "Fridge"
$ws.Range("B7").Value = "Luke Caputo G."
$ws.Range("C7").Value = "U 09"

$ws.Range("A8").Value = "Microwave"

$ws.Range("A9").Value = "Oven"
$ws.Range("B9").Value = "Danny Löser?"
$ws.Range("C9").Value = "U 08"

$ws.Range("A10").Value = "Cooking Surfaces"
$ws.Range("B10").Value = "Michael Stengel"
$ws.Range("C10").Value = "U 06"

$ws.Range("A11").Value = "Pans, Pots, etc."

$ws.Range("A12").Value = "Table & Chairs"

$ws.Range("A13").Value = "Toaster & Kettle"

$ws.Range("A14").Value = "Garbage & Area around"

$ws.Range("A15").Value = "Sink & Wall behind"

# View: zoomed in, selection on B2
$excel.ActiveWindow.Zoom = 179
$ws.Range("B2").Select()

# Page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
